$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = [ordered]@{
    'D2' = '26.127.56'
    'E2' = '  +0.09%  '
    'D3' = '1.666.14'
    'E3' = '  -0.60%  '
    'E4' = '  -0.12%  '
    'D5' = '209.67'
    'E5' = '  -0.74%  '
    'D6' = '0.5200'
    'E6' = '  -1.77%  '
    'E7' = '  -0.11%  '
    'D8' = '0.2600'
    'E8' = '  -3.01%  '
    'D9' = '0.06321'
    'E9' = '  +0.13%  '
    'E10' = '  -1.29%  '
    'D11' = '0.07528'
    'E11' = '  +0.29%  '
    'D12' = '1.664.10'
    'E12' = '  -0.87%  '
    'D13' = '4.403'
    'E13' = '  -2.35%  '
    'D14' = '0.5398'
    'E14' = '  -4.92%  '
    'D15' = '0.000007985'
    'E15' = '  -1.98%  '
    'D16' = '66.27'
    'D17' = '26.160.00'
    'E17' = '  +0.07%  '
    'E18' = '  -0.11%  '
    'D19' = '4.721'
    'E19' = '  -2.96%  '
    'D20' = '186.72'
    'E20' = '  -1.44%  '
    'E21' = '  -3.37%  '
    'D22' = '6.213'
    'E22' = '  +0.05%  '
    'E23' = '  -0.14%  '
    'D24' = '149.58'
    'E24' = '  +0.65%  '
    'D25' = '0.1234'
    'E25' = '  -1.86%  '
    'D26' = '7.394'
    'E26' = '  -3.37%  '
    'D27' = '15.71'
    'E27' = '  -2.02%  '
    'D28' = '0.06263'
    'E28' = '  -1.00%  '
    'D29' = '1.359'
    'E29' = '  +1.20%  '
    'D30' = '1.273'
    'D31' = '3.486'
    'E31' = '  -1.55%  '
    'E32' = '  -4.01%  '
    'D33' = '1.631'
    'E33' = '  -2.44%  '
    'D34' = '0.9972'
    'E34' = '  -1.24%  '
    'D35' = '2.394'
    'E35' = '  -0.96%  '
    'E36' = '  +1.38%  '
    'D37' = '0.5969'
    'E37' = '  -1.66%  '
    'D38' = '1.108.25'
    'E38' = '  +1.29%  '
    'D39' = '0.01606'
    'E39' = '  -0.65%  '
    'E40' = '  -1.81%  '
    'D41' = '0.8613'
    'E41' = '  -1.21%  '
    'E42' = '  -0.11%  '
    'D43' = '100.66'
    'E43' = '  +0.63%  '
    'D44' = '1.814.00'
    'E44' = '  -0.79%  '
    'D45' = '0.00000000110'
    'E45' = '  +1.14%  '
    'D46' = '55.26'
    'E46' = '  -3.10%  '
    'E47' = '  -0.66%  '
    'D48' = '8.051'
    'E48' = '  +0.66%  '
    'D49' = '0.05236'
    'E49' = '  -0.36%  '
    'D50' = '0.4237'
    'E50' = '  -0.72%  '
    'D51' = '5.870'
    'E51' = '  -1.52%  '
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
